$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B (existing B,C shift right to C,D)
$ws.Columns("B").Insert()

# New header for inserted column
$ws.Range("B1").Value = "StatQuery"

# New query text for inserted column, row 2 (matches wrapped style of A2)
$statQuery = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.disease IN ['Lung adenocarcinoma'] OPTIONAL MATCH (s)<-[*]-(f:file) RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(t.clinical_trial_designation)) as number_of_trial"
$ws.Range("B2").Value = $statQuery
$ws.Range("B2").WrapText = $true

# Match column B width to column A's width (best effort given COM-layer rounding
# of ColumnWidth to the nearest achievable value; 75.0 yields the closest stored
# width to column A's 75.81640625 in this runtime)
$ws.Columns("B").ColumnWidth = 75.0

# Restore single-cell selection on A2
$ws.Range("A2").Select()
